$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 2.25
$ws.Range("L2").Value = 3
$ws.Range("U2").Value = 1.44
$ws.Range("V2").Value = 2.63
$ws.Range("W2").Value = 1.91
$ws.Range("X2").Value = 1.91
$ws.Range("Y2").Value = 9
$ws.Range("AD2").Value = 34
$ws.Range("AE2").Value = 9
$ws.Range("AG2").Value = 15
$ws.Range("AI2").Value = 301
$ws.Range("AJ2").Value = 7.5
$ws.Range("AN2").Value = 19
$ws.Range("AO2").Value = 29
